# Update row values in "展览" (sheet 1) and "全部类型" (sheet 4) worksheets
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F11").Value = 55
$ws1.Range("F12").Value = 2256
$ws1.Range("F13").Value = 74

# Sheet 4: 全部类型 (All Types) - contains duplicated entries at rows 14-16
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F14").Value = 55
$ws4.Range("F15").Value = 2256
$ws4.Range("F16").Value = 74
